# Updated cryptos list on Tue Mar 26 05:07:21 UTC 2024 with GitHub Actions
# Applies latest scraped coinranking.com values to the Sheet1 data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data columns (Price, Volume(1h)) are stored as plain text in the sheet
# (e.g. "193.70", "0.120", "  +4.61%  "). Force text format first so
# Excel's COM layer doesn't "helpfully" coerce numeric-looking strings
# into numbers and strip significant trailing/leading zeros.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '70.536.08'
$ws.Range("E2").Value = '  +4.35%  '
$ws.Range("D3").Value = '3.627.94'
$ws.Range("E3").Value = '  +3.82%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = '590.88'
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("D6").Value = '193.70'
$ws.Range("E6").Value = '  +4.61%  '
$ws.Range("E7").Value = '  +1.75%  '
$ws.Range("D8").Value = '3.622.48'
$ws.Range("E8").Value = '  +3.75%  '
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").Value = '  +2.38%  '
$ws.Range("D11").Value = '0.670'
$ws.Range("E11").Value = '  +2.74%  '
$ws.Range("D12").Value = '58.33'
$ws.Range("E12").Value = '  +3.33%  '
$ws.Range("E13").Value = '  +3.72%  '
$ws.Range("D14").Value = '9.94'
$ws.Range("E14").Value = '  +4.87%  '
$ws.Range("D15").Value = '4.217.46'
$ws.Range("E15").Value = '  +4.60%  '
$ws.Range("D16").Value = '19.81'
$ws.Range("E16").Value = '  +5.49%  '
$ws.Range("D17").Value = '3.633.41'
$ws.Range("E17").Value = '  +4.19%  '
$ws.Range("D18").Value = '70.529.62'
$ws.Range("E18").Value = '  +4.51%  '
$ws.Range("D19").Value = '12.71'
$ws.Range("E19").Value = '  +4.27%  '
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("E21").Value = '  +4.13%  '
$ws.Range("D22").Value = '489.08'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '19.32'
$ws.Range("E23").Value = '  +14.04%  '
$ws.Range("D24").Value = '5.40'
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").Value = '4.46'
$ws.Range("E25").Value = '  -0.16%  '
$ws.Range("D26").Value = '91.08'
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("E27").Value = '  +6.13%  '
$ws.Range("E28").Value = '  +2.96%  '
$ws.Range("D29").Value = '9.68'
$ws.Range("E29").Value = '  +5.52%  '
$ws.Range("D30").Value = '33.10'
$ws.Range("E30").Value = '  +4.48%  '
$ws.Range("D31").Value = '7.85'
$ws.Range("E31").Value = '  +9.05%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '626.82'
$ws.Range("E32").Value = '  +5.11%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.120'
$ws.Range("E33").Value = '  +7.51%  '
$ws.Range("B34").Value = 'Cosmos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D34").Value = '12.30'
$ws.Range("E34").Value = '  +4.54%  '
$ws.Range("D35").Value = '66.09'
$ws.Range("E35").Value = '  +2.39%  '
$ws.Range("D36").Value = '39.94'
$ws.Range("E36").Value = '  +8.87%  '
$ws.Range("D37").Value = '0.414'
$ws.Range("E37").Value = '  +6.58%  '
$ws.Range("E38").Value = '  +6.42%  '
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '1.00'
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").Value = '0.146'
$ws.Range("E40").Value = '  -2.20%  '
$ws.Range("E41").Value = '  +1.00%  '
$ws.Range("D42").Value = '3.298.65'
$ws.Range("E42").Value = '  +1.21%  '
$ws.Range("D43").Value = '3.17'
$ws.Range("E43").Value = '  +8.60%  '
$ws.Range("E44").Value = '  +9.96%  '
$ws.Range("E45").Value = '  +5.22%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '3.32'
$ws.Range("E46").Value = '  +1.76%  '
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '0.139'
$ws.Range("E47").Value = '  +2.57%  '
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '2.79'
$ws.Range("E48").Value = '  +1.16%  '
$ws.Range("B49").Value = 'THORChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D49").Value = '9.18'
$ws.Range("E49").Value = '  +4.85%  '
$ws.Range("D50").Value = '3.36'
$ws.Range("E50").Value = '  +3.55%  '
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.10%  '
